$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 84, shifting existing rows (84..124) down to (85..125)
$ws.Rows.Item(84).Insert()

# Copy the style of the date cell (D85, formerly D84) onto the new D84 so the
# date format carries through to the freshly inserted row.
$ws.Range("D85").Copy()
$ws.Range("D84").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 84 with this week's data point
$ws.Cells.Item(84, 1).Value = 11
$ws.Cells.Item(84, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(84, 3).Value = "Bíobío"
$ws.Cells.Item(84, 4).Value = 45119
$ws.Cells.Item(84, 5).Value = 8
$ws.Cells.Item(84, 6).Value = 100112037
$ws.Cells.Item(84, 7).Value = "Cebollín"
$ws.Cells.Item(84, 8).Value = "Sin especificar"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 100
$ws.Cells.Item(84, 11).Value = 4500
$ws.Cells.Item(84, 12).Value = 5000
$ws.Cells.Item(84, 13).Value = 4750
$ws.Cells.Item(84, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(84, 15).Value = "Región Metropolitana"
$ws.Cells.Item(84, 16).Value = 132
$ws.Cells.Item(84, 17).Value = 36
$ws.Cells.Item(84, 18).Value = "Hortaliza"
